$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A / B (news excerpt + tag) rows 1-24
$ws.Range("A1").Value = "sa Hunyo 30"
$ws.Range("B1").Value = "DATE"
$ws.Range("A2").Value = "sa Marso 2003"
$ws.Range("B2").Value = "DATE"
$ws.Range("A3").Value = "Senado"
$ws.Range("B3").Value = "LOCATION"
$ws.Range("A4").Value = "Ayon"
$ws.Range("B4").Value = "PERSON"
$ws.Range("A5").Value = "Erap"
$ws.Range("B5").Value = "PERSON"
$ws.Range("A6").Value = "Estrada"
$ws.Range("B6").Value = "PERSON"
$ws.Range("A7").Value = "Idinagdag Villaignacio"
$ws.Range("B7").Value = "PERSON"
$ws.Range("A8").Value = "Ito"
$ws.Range("B8").Value = "PERSON"
$ws.Range("A9").Value = "Jinggoy"
$ws.Range("B9").Value = "PERSON"
$ws.Range("A10").Value = "Kongreso"
$ws.Range("B10").Value = "PERSON"
$ws.Range("A11").Value = "Magugunitang"
$ws.Range("B11").Value = "PERSON"
$ws.Range("A12").Value = "Once the SC"
$ws.Range("B12").Value = "PERSON"
$ws.Range("A13").Value = "P500,000"
$ws.Range("B13").Value = "PERSON"
$ws.Range("A14").Value = "Samantala,"
$ws.Range("B14").Value = "PERSON"
$ws.Range("A15").Value = "Sandiganbayan Special Division"
$ws.Range("B15").Value = "PERSON"
$ws.Range("A16").Value = "Senator-elect Jose Jinggoy`" Estrada`""
$ws.Range("B16").Value = "PERSON"
$ws.Range("A17").Value = "Sinabi Villaignacio"
$ws.Range("B17").Value = "PERSON"
$ws.Range("A18").Value = "Siniguro"
$ws.Range("B18").Value = "PERSON"
$ws.Range("A19").Value = "Special Division"
$ws.Range("B19").Value = "PERSON"
$ws.Range("A20").Value = "Special Prosecutor Dennis Villaignacio"
$ws.Range("B20").Value = "PERSON"
$ws.Range("A21").Value = "Ulat Malou Rongalerios"
$ws.Range("B21").Value = "PERSON"
$ws.Range("A22").Value = "Villaignacio"
$ws.Range("B22").Value = "PERSON"
$ws.Range("A23").Value = "Walang"
$ws.Range("B23").Value = "PERSON"
$ws.Range("A24").Value = "dating Pangulong Joseph Estrada"
$ws.Range("B24").Value = "PERSON"

# Column D (entity list) / F (date list) - re-ordered values
$ws.Range("D4").Value = "batang Estrada"
$ws.Range("F4").Value = "Marso 2003"
$ws.Range("D5").Value = "dating Pangulong Joseph Estrada"
$ws.Range("F5").Value = "Hunyo 30"
$ws.Range("D6").Value = "Erap"
$ws.Range("D7").Value = "Jinggoy"
$ws.Range("D8").Value = "Malou Rongalerios"
$ws.Range("D9").Value = "Senator-elect Jose `"Jinggoy`" Estrada"
$ws.Range("D10").Value = "Special Prosecutor Dennis Villaignacio"
$ws.Range("D11").Value = "Villaignacio"

# Update the used range selection / dimension to match new extent
$ws.Range("A1:B24").Select()
